# tests: added more CPU
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("opcodes")

# Mark opcodes 0x06, 0x16, 0x26, 0x36 (column H, rows 2-5) as "Tested" ("o")
# instead of "Implemented and needs integration or mocking unit tests" ("i").
$rng = $ws.Range("H2:H5")
$rng.Value = "o"

# Reflect the resulting selection state (selection H2:H5, as after editing
# the range top-to-bottom).
$ws.Range("H2:H5").Select() | Out-Null
